# Update the cryptos worksheet with the latest scraped price / volume(1h) figures,
# and fix the ordering of VeChain / TrustWalletToken rows (39-40) to match the
# refreshed coinranking.com export (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.678.33"
$ws.Range("E2").Value = "  -8.27%  "
$ws.Range("D3").Value = "1.653.29"
$ws.Range("E3").Value = "  -9.20%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").Value = "'220.42"
$ws.Range("E5").Value = "  -5.47%  "
$ws.Range("D6").Value = "'0.5091"
$ws.Range("E6").Value = "  -13.89%  "
$ws.Range("D7").Value = "'1.009"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").Value = "'0.2535"
$ws.Range("E8").Value = "  -7.54%  "
$ws.Range("D9").Value = "'21.71"
$ws.Range("E9").Value = "  -5.54%  "
$ws.Range("D10").Value = "'0.06125"
$ws.Range("E10").Value = "  -9.82%  "
$ws.Range("D11").Value = "'0.07365"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").Value = "1.640.89"
$ws.Range("E12").Value = "  -10.43%  "
$ws.Range("D13").Value = "'4.457"
$ws.Range("E13").Value = "  -4.68%  "
$ws.Range("D14").Value = "'0.5730"
$ws.Range("E14").Value = "  -8.13%  "
$ws.Range("D15").Value = "1.874.66"
$ws.Range("E15").Value = "  -9.18%  "
$ws.Range("D16").Value = "'0.000008088"
$ws.Range("E16").Value = "  -14.05%  "
$ws.Range("D17").Value = "'64.63"
$ws.Range("E17").Value = "  -13.34%  "
$ws.Range("D18").Value = "26.659.69"
$ws.Range("E18").Value = "  -7.46%  "
$ws.Range("D19").Value = "'4.963"
$ws.Range("E19").Value = "  -8.56%  "
$ws.Range("D20").Value = "'1.012"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("E21").Value = "  -7.37%  "
$ws.Range("D22").Value = "'181.10"
$ws.Range("E22").Value = "  -12.96%  "
$ws.Range("D23").Value = "'1.010"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").Value = "'6.195"
$ws.Range("E24").Value = "  -8.56%  "
$ws.Range("D25").Value = "'143.13"
$ws.Range("E25").Value = "  -7.14%  "
$ws.Range("D26").Value = "'7.600"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("E27").Value = "  -10.33%  "
$ws.Range("D28").Value = "'15.01"
$ws.Range("E28").Value = "  -7.75%  "
$ws.Range("D29").Value = "'1.329"
$ws.Range("E29").Value = "  -5.82%  "
$ws.Range("D30").Value = "'0.05809"
$ws.Range("E30").Value = "  -10.33%  "
$ws.Range("D31").Value = "'1.340"
$ws.Range("E31").Value = "  -6.55%  "
$ws.Range("D32").Value = "'3.425"
$ws.Range("E32").Value = "  -7.74%  "
$ws.Range("D33").Value = "'3.420"
$ws.Range("E33").Value = "  -7.07%  "
$ws.Range("D34").Value = "'1.573"
$ws.Range("E34").Value = "  -6.28%  "
$ws.Range("D35").Value = "'0.9813"
$ws.Range("D36").Value = "'2.429"
$ws.Range("E36").Value = "  -3.99%  "
$ws.Range("D37").Value = "'0.5971"
$ws.Range("E37").Value = "  -5.48%  "
$ws.Range("D38").Value = "'2.632"
$ws.Range("E38").Value = "  -4.20%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01576"
$ws.Range("E39").Value = "  -7.63%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.8661"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("D41").Value = "1.068.68"
$ws.Range("E41").Value = "  -5.71%  "
$ws.Range("D42").Value = "'1.012"
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("D43").Value = "'5.716"
$ws.Range("E43").Value = "  -11.23%  "
$ws.Range("D44").Value = "'95.99"
$ws.Range("E44").Value = "  -3.63%  "
$ws.Range("D45").Value = "1.784.39"
$ws.Range("E45").Value = "  -9.71%  "
$ws.Range("D46").Value = "'0.00000000108"
$ws.Range("E46").Value = "  -4.11%  "
$ws.Range("D47").Value = "'1.014"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").Value = "'55.16"
$ws.Range("E48").Value = "  -8.41%  "
$ws.Range("D49").Value = "'0.4383"
$ws.Range("E49").Value = "  -2.93%  "
$ws.Range("D50").Value = "'0.05205"
$ws.Range("E50").Value = "  -4.92%  "
$ws.Range("D51").Value = "'7.764"
$ws.Range("E51").Value = "  -6.08%  "
